# edit.ps1
# Applies three changes to Roadmap.docx:
#   1. Remove the stray "_GoBack" bookmark in the empty paragraph after "Diogo | Arnaldo".
#   2. Split "...consecutivos)" into "...consecutivos" + "=" + (bookmark "_GoBack") + ")".
#   3. Merge the "- " + "S" + "ó conjunto de eventos" runs into a single "- Só conjunto de eventos" run.
#
# NOTE on technique: this engine re-coalesces all adjacent runs that share identical
# formatting (rPr) within a paragraph whenever a run boundary is touched by an
# insert/delete. So the reliable way to end up with an *exact* target run layout is:
#   (a) let the paragraph coalesce into the minimum number of runs however it likes,
#   (b) then force precise split points back in with a harmless formatting toggle
#       (Bold On then Off) applied to exactly the sub-range that should become its
#       own run. Toggling Bold like this leaves no residue in rPr (unlike e.g. Italic).

$d = $word.ActiveDocument

function Force-Split([object]$range) {
    # Forces run boundaries at both ends of $range without altering visible formatting.
    $range.Bold = 1
    $range.Bold = 0
}

# ---------------------------------------------------------------------------
# Change 1: remove the _GoBack bookmark that sits in the empty paragraph right
# after "Diogo | Arnaldo".
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# Change 2: "... que deteta que o jogador marcou golo em 2 jogos consecutivos)"
# becomes three runs: "...consecutivos" / "=" / ")" with a (new) _GoBack
# bookmark sitting between the "=" run and the ")" run.
# ---------------------------------------------------------------------------
$search = $d.Content
$found = $search.Find.Execute("consecutivos)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'consecutivos)' text"
}

$pRange = $search.Paragraphs(1).Range
$pStart = $pRange.Start
$pEnd = $pRange.End

# Position right before the final ')' character (pEnd is just past the
# paragraph mark, so the ')' character itself sits at pEnd-2).
$parenPos = $pEnd - 2
$insertion = $d.Range($parenPos, $parenPos)
$insertion.InsertAfter("=")

# The paragraph text (ignoring the trailing paragraph mark) now reads:
#   "- Captura 100% funcional de pelo menos 1 tipo de badge (para efeitos
#    demonstrativos, vamos ter apenas um badge que deteta que o jogador
#    marcou golo em 2 jogos consecutivos=)"
# which the engine will have coalesced into a single run covering the whole
# paragraph. Re-split it back into the desired run layout using offsets
# relative to the (fresh) paragraph start.
$offsets = @(0, 49, 54, 105, 110, 171, 172, 173)
for ($i = 0; $i -lt $offsets.Length - 1; $i++) {
    $segStart = $pStart + $offsets[$i]
    $segEnd = $pStart + $offsets[$i + 1]
    Force-Split ($d.Range($segStart, $segEnd))
}

# Insert the _GoBack bookmark between the "=" run (offsets 171-172) and the
# ")" run (offsets 172-173).
$bmPos = $pStart + 172
$d.Bookmarks.Add("_GoBack", $d.Range($bmPos, $bmPos))

# ---------------------------------------------------------------------------
# Change 3: merge "- " + "S" + "ó conjunto de eventos" into a single run
# "- Só conjunto de eventos", while keeping the preceding "Notificações " run
# separate.
# ---------------------------------------------------------------------------
$search2 = $d.Content
$found2 = $search2.Find.Execute("conjunto de eventos", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find 'conjunto de eventos' text"
}

$pRange2 = $search2.Paragraphs(1).Range
$pStart2 = $pRange2.Start

# "Notificações " is 13 characters long; touch the run boundary right after it
# (inserting then immediately removing a throwaway character) to force the
# engine to re-coalesce the remaining identically-formatted runs into one.
$boundary = $pStart2 + 13
$tmp = $d.Range($boundary, $boundary)
$tmp.InsertBefore("Z")
$tmpChar = $d.Range($boundary, $boundary + 1)
$tmpChar.Text = ""

# Now the whole paragraph is a single run; split "Notificações " back off so
# it remains its own run, leaving "- Só conjunto de eventos" as the other.
Force-Split ($d.Range($pStart2, $pStart2 + 13))

Write-Output "Done"
